# This workbook (FuzzyMatchingOut.xlsx) was originally produced by a
# non-Excel tool (openpyxl) and this change represents the file being
# opened in Excel, having its columns auto-sized to fit their contents,
# and then saved again (the active selection ends up on G20).
#
# Recreate that by auto-fitting the data columns (B:E) and leaving the
# selection on cell G20 before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-size the columns that hold data (B=RegionAbbrev/Region, C/D = grant
# text, E = match ratio) so their widths reflect their content, the same
# way Excel does automatically when a worksheet produced by another tool
# is first opened and saved.
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 42.166666666666664
$ws.Columns.Item(4).ColumnWidth = 43
$ws.Columns.Item(5).ColumnWidth = 24.5

# Leave the selection where the author last left it when the file was
# saved.
$ws.Range("G20").Select() | Out-Null
